$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (shifting existing row 5 "population" down to row 6,
# and Oporavljeni/Testirani/Smrtni sl. down to rows 7/8/9)
$ws.Rows.Item(5).Insert()

# New row 5: stringency_index
$ws.Range("A5").Value = "stringency_index"
$ws.Range("B5").Value = 308
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

# Update B column for date, total_cases, new_cases (rows 2-4) -> 308
$ws.Range("B2").Value = 308
$ws.Range("B3").Value = 308
$ws.Range("B4").Value = 308

# Row 6: population -> B=308
$ws.Range("B6").Value = 308

# Row 7: Oporavljeni
$ws.Range("C7").Value = 77
$ws.Range("D7").Value = 0.3333333333333333

# Row 8: Testirani
$ws.Range("C8").Value = 77
$ws.Range("D8").Value = 0.3333333333333333

# Row 9: Smrtni sl.
$ws.Range("C9").Value = 76
$ws.Range("D9").Value = 0.3275862068965517
